$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Every cell across the workbook that showed this status moves together,
#    since it was a single shared string being edited in place.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Report generated for handback: populate "Latest Target File",
#    "Latest Handback File" and "Latest Handback DateTime" on both the
#    zh-cn and de-de detail sheets for the two rows.
# ---------------------------------------------------------------------------

$md1 = "88552b79-3232-45b1-9dbb-f65282794b4c.md"
$md2 = "ce9e5566-f3ff-4db1-9103-3cc2fc4ec57b.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/c22f77d18ed4bb10a8e3e44df5cbf9a9585d1b9e/e2e/88552b79-3232-45b1-9dbb-f65282794b4c.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/c22f77d18ed4bb10a8e3e44df5cbf9a9585d1b9e/e2e/ce9e5566-f3ff-4db1-9103-3cc2fc4ec57b.md"

# zh-cn sheet
$wsZh.Range("I2").Value = "88552b79-3232-45b1-9dbb-f65282794b4c.4e69898a37946567d0d9bc9b650b911e5a646c1f.zh-cn.xlf"
$wsZh.Range("J2").Value = "2016-07-26 08:23:14"
$wsZh.Hyperlinks.Add($wsZh.Range("H2"), $url1, "", "", $md1)

$wsZh.Range("I3").Value = "ce9e5566-f3ff-4db1-9103-3cc2fc4ec57b.d93bf6ab61dcd3b74c5fe2538052fd1321ecf7ac.zh-cn.xlf"
$wsZh.Range("J3").Value = "2016-07-26 08:23:14"
$wsZh.Hyperlinks.Add($wsZh.Range("H3"), $url2, "", "", $md2)

# de-de sheet
$wsDe.Range("I2").Value = "88552b79-3232-45b1-9dbb-f65282794b4c.4e69898a37946567d0d9bc9b650b911e5a646c1f.de-de.xlf"
$wsDe.Range("J2").Value = "2016-07-26 08:23:32"
$wsDe.Hyperlinks.Add($wsDe.Range("H2"), $url1, "", "", $md1)

$wsDe.Range("I3").Value = "ce9e5566-f3ff-4db1-9103-3cc2fc4ec57b.d93bf6ab61dcd3b74c5fe2538052fd1321ecf7ac.de-de.xlf"
$wsDe.Range("J3").Value = "2016-07-26 08:23:32"
$wsDe.Hyperlinks.Add($wsDe.Range("H3"), $url2, "", "", $md2)

# ---------------------------------------------------------------------------
# 3. Widen columns that now hold the longer status text / file names.
# ---------------------------------------------------------------------------

# Overview: zh-cn (E) and de-de (F) status columns
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# zh-cn: Status (C), Latest Target File (H), Latest Handback File (I)
$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(8).ColumnWidth = 39.14
$wsZh.Columns.Item(9).ColumnWidth = 39.14

# de-de: Status (C), Latest Target File (H), Latest Handback File (I)
$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(8).ColumnWidth = 39.14
$wsDe.Columns.Item(9).ColumnWidth = 39.14
